# ESPUMOSO.xlsx update
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Delete the "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

# Remember which sheet should stay active/selected (first sheet, "Paineis DARQ").
$mainSheet = $wb.Worksheets.Item("Paineis DARQ")

# Remove the sheet that is no longer needed.
[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Rename the remaining sheets.
$mainSheet.Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the original active sheet selected (deleting a sheet can shift the
# active tab in some hosts), matching the workbook's original view state.
[void]$mainSheet.Activate()
